# Adaptive-Sonar-Route-Planning "antwerpen" decomposition fix-up.
#
# The sheet holds a route (columns A index, B x-coordinate, C y-coordinate).
# This edit adds helper columns:
#   G1 = MIN(B1:B50)   H1 = MIN(C1:C50)
#   D1:D50 = array formula  B1:B50-G1   (x relative to the minimum x)
#   E1:E50 = array formula  C1:C50-H1   (y relative to the minimum y)
# replacing the old "#" placeholder text that used to sit in D/E.
# A handful of rows (9, 11-46) keep a leftover "#" marker, now in column F,
# from the manual row-by-row clean-up described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear any pre-existing formatting on D:E first -- leftover style bookkeeping
# from the old "#" placeholder cells otherwise bleeds back in once a
# NumberFormat is touched anywhere in the column.
$ws.Range("D1:E50").ClearFormats()

# --- helper columns with the running minimums ---
$ws.Range("G1").Formula = "=MIN(B1:B50)"
$ws.Range("H1").Formula = "=MIN(C1:C50)"
$ws.Range("G1").NumberFormat = "0"
$ws.Range("H1").NumberFormat = "0"

# --- legacy (CSE-style) array formulas spilling down D1:D50 / E1:E50 ---
$ws.Range("D1:D50").FormulaArray = "=B1:B50-G1"
$ws.Range("E1:E50").FormulaArray = "=C1:C50-H1"

# A few cells close to the top retain the "0" integer number format that was
# already used elsewhere on the sheet (style index reused automatically).
$ws.Range("D1").NumberFormat = "0"
$ws.Range("E1").NumberFormat = "0"
$ws.Range("D2").NumberFormat = "0"
$ws.Range("D3").NumberFormat = "0"

# --- leftover "#" markers in column F for the rows that were hand-patched ---
$markerRows = @(9,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46)
foreach ($r in $markerRows) {
    $ws.Range("F$r").Value = "#"
}

# --- a few stray styled-but-empty cells further down the sheet ---
$ws.Range("F54").NumberFormat = "0"
$ws.Range("F55").NumberFormat = "0"
$ws.Range("F57").NumberFormat = "0"
$ws.Range("F58").NumberFormat = "0"

# --- restore the view: scroll back to the top and select H12 ---
$ws.Range("H12").Select()
